$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.172.18"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "1.783.69"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3956"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3425"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.84"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("E10").Value = "  -3.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07445"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.441"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.79%  "

$ws.Range("D15").Value = "1.780.34"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.087"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("E17").Value = "  -3.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06671"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.14%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.492"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("D23").Value = "27.178.05"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("E24").Value = "  -6.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.375"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.47%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.499"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.454"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "

$ws.Range("D30").Value = "1.981.58"
$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.970"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("E33").Value = "  -6.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08752"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.610"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.69%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.398"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.72%  "

$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02378"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6788"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06370"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("E41").Value = "  -2.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.239"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.431"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.91%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6384"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.69%  "

$ws.Range("E47").Value = "  -2.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.133"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.16%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07105"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.27%  "
